$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lgals1"
$ws.Range("C2").Value = "Cd69"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 37.833119
$ws.Range("H2").Value = 113.499357
$ws.Range("I2").Value = 0.2771305381131279
$ws.Range("J2").Value = 0.2771305381131279
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 36.51526533333334
$ws.Range("N2").Value = 109.545796
$ws.Range("O2").Value = 0.9971632597677756
$ws.Range("P2").Value = 0.9971632597677755
$ws.Range("Q2").Value = 1381.486378672575
$ws.Range("R2").Value = 12433.37740805317
$ws.Range("S2").Value = 0.2763443907660844
$ws.Range("T2").Value = 0.2763443907660844

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lgals1"
$ws.Range("C3").Value = "Cd69"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 37.833119
$ws.Range("H3").Value = 113.499357
$ws.Range("I3").Value = 0.2771305381131279
$ws.Range("J3").Value = 0.2771305381131279
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.103879
$ws.Range("N3").Value = 0.311637
$ws.Range("O3").Value = 0.002836740232224432
$ws.Range("P3").Value = 0.002836740232224432
$ws.Range("Q3").Value = 3.930066568601
$ws.Range("R3").Value = 35.370599117409
$ws.Range("S3").Value = 0.0007861473470435163
$ws.Range("T3").Value = 0.0007861473470435163

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Lgals1"
$ws.Range("C4").Value = "Cd69"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 63.88336466666667
$ws.Range("H4").Value = 191.650094
$ws.Range("I4").Value = 0.4679506129682439
$ws.Range("J4").Value = 0.467950612968244
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 36.51526533333334
$ws.Range("N4").Value = 109.545796
$ws.Range("O4").Value = 0.9971632597677756
$ws.Range("P4").Value = 0.9971632597677755
$ws.Range("Q4").Value = 2332.718011189425
$ws.Range("R4").Value = 20994.46210070483
$ws.Range("S4").Value = 0.4666231586377428
$ws.Range("T4").Value = 0.4666231586377428

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Lgals1"
$ws.Range("C5").Value = "Cd69"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.88336466666667
$ws.Range("H5").Value = 191.650094
$ws.Range("I5").Value = 0.4679506129682439
$ws.Range("J5").Value = 0.467950612968244
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.103879
$ws.Range("N5").Value = 0.311637
$ws.Range("O5").Value = 0.002836740232224432
$ws.Range("P5").Value = 0.002836740232224432
$ws.Range("Q5").Value = 6.636140038208667
$ws.Range("R5").Value = 59.72526034387801
$ws.Range("S5").Value = 0.001327454330501102
$ws.Range("T5").Value = 0.001327454330501102

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Lgals1"
$ws.Range("C6").Value = "Cd69"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 34.80083866666666
$ws.Range("H6").Value = 104.402516
$ws.Range("I6").Value = 0.2549188489186281
$ws.Range("J6").Value = 0.2549188489186282
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 36.51526533333334
$ws.Range("N6").Value = 109.545796
$ws.Range("O6").Value = 0.9971632597677756
$ws.Range("P6").Value = 0.9971632597677755
$ws.Range("Q6").Value = 1270.76185773586
$ws.Range("R6").Value = 11436.85671962274
$ws.Range("S6").Value = 0.2541957103639483
$ws.Range("T6").Value = 0.2541957103639484

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Lgals1"
$ws.Range("C7").Value = "Cd69"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 34.80083866666666
$ws.Range("H7").Value = 104.402516
$ws.Range("I7").Value = 0.2549188489186281
$ws.Range("J7").Value = 0.2549188489186282
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.103879
$ws.Range("N7").Value = 0.311637
$ws.Range("O7").Value = 0.002836740232224432
$ws.Range("P7").Value = 0.002836740232224432
$ws.Range("Q7").Value = 3.615076319854666
$ws.Range("R7").Value = 32.535686878692
$ws.Range("S7").Value = 0.000723138554679814
$ws.Range("T7").Value = 0.0007231385546798142
